# Adicionando os loops e array de ambientes
#
# Atualiza o nome do ambiente cadastrado (celula A2) e remove o sublinhado
# do estilo aplicado a essa celula.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Troca o texto do ambiente cadastrado.
$ws.Range("A2").Value = "Teste01"

# O estilo da celula A2 usava uma fonte sublinhada; remove o sublinhado
# (fontId passa de 2 -> 1 no styles.xml).
$ws.Range("A2").Font.Underline = -4142
